$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-15) for columns D, L, M, N, O, P, Q, R, S, T
$data = @(
    @{ Row=2;  D=44319; L='Especial'; M=120; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=3;  D=44294; L='Primera';  M=50;  N=12000; O=12000; P=12000; Q='$/caja 15 kilos granel'; R='Región Metropolitana'; S=800;  T=15 },
    @{ Row=4;  D=44348; L='Especial'; M=200; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=5;  D=44291; L='Primera';  M=150; N=12000; O=12000; P=12000; Q='$/caja 15 kilos granel'; R='Región Metropolitana'; S=800;  T=15 },
    @{ Row=6;  D=44328; L='Especial'; M=250; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=7;  D=44326; L='Especial'; M=300; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=8;  D=44354; L='Primera';  M=100; N=18000; O=18000; P=18000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1000; T=18 },
    @{ Row=9;  D=44355; L='Especial'; M=50;  N=18000; O=18000; P=18000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1000; T=18 },
    @{ Row=10; D=44340; L='Primera';  M=230; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=11; D=44299; L='Primera';  M=100; N=15000; O=15000; P=15000; Q='$/caja 15 kilos granel'; R='Provincia de Curicó'; S=1000; T=15 },
    @{ Row=12; D=44316; L='Especial'; M=300; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 },
    @{ Row=13; D=44358; L='Especial'; M=150; N=18000; O=18000; P=18000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1000; T=18 },
    @{ Row=14; D=44358; L='Primera';  M=100; N=17000; O=17000; P=17000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=944;  T=18 },
    @{ Row=15; D=44342; L='Especial'; M=300; N=20000; O=20000; P=20000; Q='$/caja 18 kilos granel'; R='Provincia de Limarí'; S=1111; T=18 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("L$r").Value = $entry.L
    $ws.Range("M$r").Value = $entry.M
    $ws.Range("N$r").Value = $entry.N
    $ws.Range("O$r").Value = $entry.O
    $ws.Range("P$r").Value = $entry.P
    $ws.Range("Q$r").Value = $entry.Q
    $ws.Range("R$r").Value = $entry.R
    $ws.Range("S$r").Value = $entry.S
    $ws.Range("T$r").Value = $entry.T
}
